# Update countries & provincias Spain
# Refresh the COVID-19 country statistics with the 20:52 data pull.
# - Updates the "last updated" timestamp
# - Updates numeric figures for several countries whose totals changed
# - Three countries (Burkina Faso, Somalia, Suazilandia) overtook their
#   neighbour in the case-count ranking and so swap places with them

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CountryRow([int]$row, [string]$country, [double]$casosTotales, [double]$nuevosCasos, [double]$casosActivos, [double]$recuperados, [double]$casosCriticos, [double]$muertesHoy, [double]$muertes) {
    $data = New-Object 'object[,]' 1,8
    $data[0,0] = $country
    $data[0,1] = $casosTotales
    $data[0,2] = $nuevosCasos
    $data[0,3] = $casosActivos
    $data[0,4] = $recuperados
    $data[0,5] = $casosCriticos
    $data[0,6] = $muertesHoy
    $data[0,7] = $muertes
    $ws.Range($ws.Cells.Item($row, 1), $ws.Cells.Item($row, 8)).Value = $data
}

# Timestamp banner (row 1)
$ws.Range("A1").Value = "Datos actualizados a 21 de Abril de 2020 a las 20:52"

# Straightforward numeric refreshes (no re-ordering)
Set-CountryRow 4  "Estados Unidos" 804759 12000 81810 678954 14016 1481 43995
Set-CountryRow 22 "Irlanda"        16040  388   9233  6077   315   43   730
Set-CountryRow 56 "Marruecos"      3209   163   393   2671   1     2    145
Set-CountryRow 65 "Barein"         1973   66    784   1182   2     0    7
Set-CountryRow 68 "Uzbekistan"     1678   51    357   1315   8     1    6

# Burkina Faso overtakes Bolivia and Kirguistan
Set-CountryRow 99  "Burkina Faso" 600 19 362 200 0 0 38
Set-CountryRow 100 "Bolivia"      598 34 37  527 3 1 34
Set-CountryRow 101 "Kirguistan"   590 22 216 367 5 0 7

# Somalia overtakes Venezuela, Vietnam, Mali and Tanzania
Set-CountryRow 120 "Somalia"   286 49 4   274 2 0 8
Set-CountryRow 121 "Venezuela" 285 0  117 158 4 0 10
Set-CountryRow 122 "Vietnam"   268 0  216 52  8 0 0
Set-CountryRow 123 "Mali"      258 12 57  187 0 0 14
Set-CountryRow 124 "Tanzania"  254 0  11  233 4 0 10

# Suazilandia overtakes Zimbabue and Angola
Set-CountryRow 175 "Suazilandia" 31 7 8 22 0 0 1
Set-CountryRow 176 "Zimbabue"    25 0 2 20 0 0 3
Set-CountryRow 177 "Angola"      24 0 6 16 0 0 2
